$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "Training Dashboard" - refresh the "PERIOD TO EXPIRE" (H) and
# "LAST UPDATE" (I) columns for rows 3-14 (the LAST UPDATE date moved from
# 08-Sep-2025 to 16-Sep-2025, eight days later, so PERIOD TO EXPIRE dropped
# by 8 for every row).
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Training Dashboard")

$periodUpdates = @{
    3  = 426
    4  = 408
    5  = 219
    6  = 400
    7  = 210
    8  = 349
    9  = 359
    10 = 395
    11 = -103
    12 = -343
    13 = -41
    14 = -41
}

foreach ($row in $periodUpdates.Keys) {
    $ws1.Range("H$row").Value = $periodUpdates[$row]
    $ws1.Range("I$row").Value = "16-Sep-2025"
    $ws1.Range("I$row").NumberFormat = "@"
}

# ---------------------------------------------------------------------------
# Sheet 2: "Exam Dashboard" - the exams are no longer outdated, update the
# comments and narrow the COMMENTS column now that the text is shorter.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Exam Dashboard")

$ws2.Range("E3").Value = "date is valid"
$ws2.Range("E4").Value = "date is valid"

# NOTE: Excel's ColumnWidth (character units) is offset from the raw OOXML
# <col width> by the default font's padding (5/6 of a character here), so
# asking for a stored width of 15 means setting ColumnWidth to 15 - 5/6.
$ws2.Columns.Item(5).ColumnWidth = 14.166666666666666

# ---------------------------------------------------------------------------
# Header/title styling: the workbook now shares a single bold font between
# the title row and the header row, so both render in white (size 11)
# instead of the title's previous black, size-14 look.
# ---------------------------------------------------------------------------
foreach ($ws in $wb.Worksheets) {
    $ws.Range("A1").Font.Color = 16777215
    $ws.Range("A1").Font.Bold = $true
    $ws.Range("A1").Font.Size = 11

    $ws.Range("A2:K2").Font.Color = 16777215
    $ws.Range("A2:K2").Font.Bold = $true
    $ws.Range("A2:K2").Font.Size = 11
}
